$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 821.9375
$ws.Range("I129").Value = 385.22223
$ws.Range("J129").Value = 1383.4286
$ws.Range("K129").Value = 1155.66669
$ws.Range("L129").Value = 4150.2858
$ws.Range("M129").Value = 3844.33331
$ws.Range("N129").Value = -14150.2858
$ws.Range("H132").Value = 1941.8649
$ws.Range("I132").Value = 1000.76
$ws.Range("K132").Value = 3002.28
$ws.Range("M132").Value = -472.2799999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 151
$ws.Range("I5").Value = 35.285713
$ws.Range("J5").Value = 556
$ws.Range("K5").Value = 35.285713
$ws.Range("L5").Value = 556
$ws.Range("M5").Value = 76.714287
$ws.Range("N5").Value = -780
$ws.Range("H74").Value = 1232.5834
$ws.Range("I74").Value = 939.2941
$ws.Range("J74").Value = 1944.8572
$ws.Range("K74").Value = 939.2941
$ws.Range("L74").Value = 1944.8572
$ws.Range("M74").Value = -65.29409999999996
$ws.Range("N74").Value = -3692.8572
$ws.Range("H77").Value = 1232.5834
$ws.Range("I77").Value = 939.2941
$ws.Range("J77").Value = 1944.8572
$ws.Range("K77").Value = 4696.470499999999
$ws.Range("L77").Value = 9724.286
$ws.Range("M77").Value = -328.4704999999994
$ws.Range("N77").Value = -18460.286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 151
$ws.Range("I4").Value = 35.285713
$ws.Range("J4").Value = 556
$ws.Range("K4").Value = 35.285713
$ws.Range("L4").Value = 556
$ws.Range("M4").Value = 79.714287
$ws.Range("N4").Value = -786
$ws.Range("H20").Value = 40174.668
$ws.Range("I20").Value = 52221.3
$ws.Range("J20").Value = 5755.7144
$ws.Range("K20").Value = 52221.3
$ws.Range("L20").Value = 5755.7144
$ws.Range("M20").Value = -51974.3
$ws.Range("N20").Value = -6249.7144
$ws.Range("H22").Value = 13637.5
$ws.Range("I22").Value = 13637.5
$ws.Range("K22").Value = 13637.5
$ws.Range("M22").Value = -13464.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.42856999999999
$ws.Range("I7").Value = 66.42856999999999
$ws.Range("K7").Value = 66.42856999999999
$ws.Range("M7").Value = 46.57143000000001
$ws.Range("H22").Value = 304.35
$ws.Range("I22").Value = 245.9375
$ws.Range("J22").Value = 538
$ws.Range("K22").Value = 245.9375
$ws.Range("L22").Value = 538
$ws.Range("M22").Value = 104.0625
$ws.Range("N22").Value = -1238
$ws.Range("H31").Value = 2383.0393
$ws.Range("I31").Value = 1866.5122
$ws.Range("J31").Value = 4500.8
$ws.Range("K31").Value = 1866.5122
$ws.Range("L31").Value = 4500.8
$ws.Range("M31").Value = -1571.5122
$ws.Range("N31").Value = -5090.8
$ws.Range("H34").Value = 2383.0393
$ws.Range("I34").Value = 1866.5122
$ws.Range("J34").Value = 4500.8
$ws.Range("K34").Value = 1866.5122
$ws.Range("L34").Value = 4500.8
$ws.Range("M34").Value = -1664.5122
$ws.Range("N34").Value = -4904.8
$ws.Range("H107").Value = 409.4
$ws.Range("I107").Value = 433.44446
$ws.Range("J107").Value = 193
$ws.Range("K107").Value = 433.44446
$ws.Range("L107").Value = 193
$ws.Range("M107").Value = 1486.55554
$ws.Range("N107").Value = -4033
$ws.Range("H134").Value = 2167
$ws.Range("I134").Value = 1867.9286
$ws.Range("J134").Value = 2690.375
$ws.Range("K134").Value = 5603.7858
$ws.Range("L134").Value = 8071.125
$ws.Range("M134").Value = -3068.7858
$ws.Range("N134").Value = -13141.125
$ws.Range("H135").Value = 40778.75
$ws.Range("J135").Value = 40778.75
$ws.Range("L135").Value = 40778.75
$ws.Range("N135").Value = -50918.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 15800
$ws.Range("I70").Value = 18500
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 55500
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -55185
$ws.Range("N70").Value = -15630
$ws.Range("H73").Value = 15800
$ws.Range("I73").Value = 18500
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 55500
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -54408
$ws.Range("N73").Value = -17184
$ws.Range("H113").Value = 659.6923
$ws.Range("I113").Value = 655.5333000000001
$ws.Range("J113").Value = 665.36365
$ws.Range("K113").Value = 1966.5999
$ws.Range("L113").Value = 1996.09095
$ws.Range("M113").Value = 203.4000999999998
$ws.Range("N113").Value = -6336.09095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2181
$ws.Range("I122").Value = 2067.7
$ws.Range("J122").Value = 2407.6
$ws.Range("K122").Value = 6203.099999999999
$ws.Range("L122").Value = 7222.799999999999
$ws.Range("M122").Value = -3753.099999999999
$ws.Range("N122").Value = -12122.8
$ws.Range("H132").Value = 2457.1035
$ws.Range("I132").Value = 2018.875
$ws.Range("J132").Value = 2996.4614
$ws.Range("K132").Value = 6056.625
$ws.Range("L132").Value = 8989.3842
$ws.Range("M132").Value = -3526.625
$ws.Range("N132").Value = -14049.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3675.111
$ws.Range("I7").Value = 3565.75
$ws.Range("K7").Value = 3565.75
$ws.Range("M7").Value = -3453.75
$ws.Range("H94").Value = 17878
$ws.Range("J94").Value = 17878
$ws.Range("L94").Value = 17878
$ws.Range("N94").Value = -19230
$ws.Range("H122").Value = 5199.9414
$ws.Range("I122").Value = 4415.3076
$ws.Range("J122").Value = 7750
$ws.Range("K122").Value = 13245.9228
$ws.Range("L122").Value = 23250
$ws.Range("M122").Value = -10795.9228
$ws.Range("N122").Value = -28150
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 105638.08
$ws.Range("J125").Value = 105638.08
$ws.Range("L125").Value = 105638.08
$ws.Range("N125").Value = -115478.08
$ws.Range("H126").Value = 3675.111
$ws.Range("I126").Value = 3565.75
$ws.Range("K126").Value = 10697.25
$ws.Range("M126").Value = -8227.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 36688.5
$ws.Range("J109").Value = 36688.5
$ws.Range("L109").Value = 36688.5
$ws.Range("N109").Value = -39462.5
$ws.Range("H122").Value = 2447.2778
$ws.Range("I122").Value = 2151.1853
$ws.Range("J122").Value = 3335.5557
$ws.Range("K122").Value = 6453.5559
$ws.Range("L122").Value = 10006.6671
$ws.Range("M122").Value = -4003.5559
$ws.Range("N122").Value = -14906.6671
